# Updated symbol list (Price / Volume(1h)) per commit diff.
# Values are stored as literal text (matching the original inlineStr
# cells), so each assignment is apostrophe-prefixed to force text
# entry, then the cosmetic "number stored as text" style flag is
# cleared by resetting the cell back to the Normal style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'285.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.85%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'29.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'3.69%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.063"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.76%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06713"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.95%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.316"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.07%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.439"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.83%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.386"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-3.09%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.8998"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-2.43%"
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'3.13%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07095"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'9.20%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07642"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.59%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.02925"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'5.57%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.08990"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.20%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001582"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.03%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.04499"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.65%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.0006462"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.76%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.006192"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'2.15%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.449"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.06%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'2.230"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.37%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'1.32%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.1319"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'1.06%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'3.900"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-2.22%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.1558"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.83%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001202"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'1.34%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004369"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-1.20%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'-6.56%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0001617"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'-0.14%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.04246"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.83%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006819"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.17%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1239"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'1.10%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'2.57%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01260"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'9.81%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005753"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'2.24%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.966"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'0.02%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'15.39%"
$ws.Range("E47").Style = "Normal"
